$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.226.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.604.15"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.41"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0813"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.828.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.593.36"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.71%  "
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.198.37"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.06"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.53%  "
$ws.Range("E18").Value = "  +0.84%  "
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "200.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.59%  "
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  +2.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.62%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  -2.19%  "
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("E29").Value = "  +1.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0491"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.92%  "
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("E32").Value = "  +2.51%  "
$ws.Range("E33").Value = "  -1.25%  "
$ws.Range("E34").Value = "  +0.53%  "
$ws.Range("E35").Value = "  +1.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.160.93"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.99%  "
$ws.Range("E37").Value = "  +3.43%  "
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.32"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("E40").Value = "  +0.39%  "
$ws.Range("E41").Value = "  +0.51%  "
$ws.Range("E42").Value = "  +0.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.740.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.76"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.29%  "
$ws.Range("E46").Value = "  +16.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.54"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.35%  "
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("E51").Value = "  -0.09%  "
